$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 56
$ws.Range("F2").Value = 38
$ws.Range("H2").Value = 50

$ws.Range("E3").Value = 26

$ws.Range("E9").Value = 29
$ws.Range("F9").Value = 15
$ws.Range("H9").Value = 24

$ws.Range("E15").Value = 168
$ws.Range("F15").Value = 94
$ws.Range("H15").Value = 135

$ws.Range("E17").Value = 131
$ws.Range("F17").Value = 68
$ws.Range("H17").Value = 100

$ws.Range("E19").Value = 67

$ws.Range("E24").Value = 26

$ws.Range("E26").Value = 33

$ws.Range("E37").Value = 59
$ws.Range("F37").Value = 36
$ws.Range("H37").Value = 48

$ws.Range("E38").Value = 83
$ws.Range("F38").Value = 20
$ws.Range("H38").Value = 40

$ws.Range("E42").Value = 40

$ws.Range("F45").Value = 15
$ws.Range("H45").Value = 22

$ws.Range("F49").Value = 42
$ws.Range("H49").Value = 59

$ws.Range("E57").Value = 17

$ws.Range("E61").Value = 31

$ws.Range("E63").Value = 42
$ws.Range("F63").Value = 15
$ws.Range("H63").Value = 23

$ws.Range("E69").Value = 17

$ws.Range("E77").Value = 60

$ws.Range("E88").Value = 28
$ws.Range("F88").Value = 17
$ws.Range("H88").Value = 25

$ws.Range("E89").Value = 48

$wb.Save()
